# Update cached datetimeFigureOut field text across the deck's
# Slide Master, all Custom Layouts (slide layouts), and the Notes Master.
# ("ו'/תמוז/תשפ"ג" -> "כ'/תמוז/תשפ"ג")
$p = $ppt.ActivePresentation
$newDate = "כ'/תמוז/תשפ`"ג"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every Custom Layout hanging off the Slide Master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j).Shapes
}

# Notes Master
Update-DatePlaceholders $p.NotesMaster.Shapes

# ---------------------------------------------------------------------
# Slide 17 ("Counter objects - cont'd"): bold the update()/subtract()/
# elements() method names referenced in the bullet text.
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$shp17 = $s17.Shapes.Item(2)
$tr17 = $shp17.TextFrame.TextRange

$para7 = $tr17.Paragraphs(7, 1)
$para7.Characters(212, 9).Font.Bold = $true   # "update() "
$para7.Characters(225, 11).Font.Bold = $true  # "subtract() "

$para11 = $tr17.Paragraphs(11, 1)
$para11.Characters(5, 11).Font.Bold = $true   # "elements() "

Write-Host "Edit complete"
